$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ZM")

$ws.Range("B7").Value = 38000000.0
$ws.Range("C7").Value = 31000000.0
$ws.Range("D7").Value = 36562000.0
$ws.Range("E7").Value = 46262000.0
$ws.Range("F7").Value = 55485000.0
